$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4  - Log / Implementação dos recursos de log de exceções
$ws.Range("D4").Value = "Logs de response (erros de servidor e da API)"

# Row 11 - Segurança / Definição do modelo de autenticação
$ws.Range("D11").Value = "Autenticação de usuário"

# Row 16 - Hardware/Provisionamento Cloud
$ws.Range("D16").Value = "Azure (Aplicativos Web e Banco de Dados)"

# Row 17 - Cliente (Notebook, Celular, Browser)
$ws.Range("D17").Value = "Browser (Chrome, Firefox e Safari)"

# Row 25 - Testes / Processo e ferramenta para realização dos Testes
$ws.Range("D25").Value = "Testes unitários das principais regras de negócio"

# Reset selection back to the top-left cell of the frozen pane
$ws.Range("A1").Select()
